$wb = $excel.ActiveWorkbook

# --- raw_influent_compound_conc sheet (sheet2) ---
$wsConc = $wb.Worksheets.Item("raw_influent_compound_conc")

# Column A (Variable names) first so new shared-strings are appended in this order
$wsConc.Range("A20").Value = "inf_BOD"
$wsConc.Range("A21").Value = "inf_COD"
$wsConc.Range("A22").Value = "inf_TKN"
$wsConc.Range("A23").Value = "inf_TN"
$wsConc.Range("A24").Value = "inf_TP"
$wsConc.Range("A25").Value = "inf_TSS"
$wsConc.Range("A26").Value = "inf_VSS"

# Column B (numeric values)
$wsConc.Range("B20").Value = 379.82797943309998
$wsConc.Range("B21").Value = 1760.4
$wsConc.Range("B22").Value = 154.42099999999999
$wsConc.Range("B23").Value = 174.42099999999999
$wsConc.Range("B24").Value = 166.740350277819
$wsConc.Range("B25").Value = 1182.6926908329499
$wsConc.Range("B26").Value = 771.45741272220005

# Column C (descriptions) last so they become shared strings 83-89
$wsConc.Range("C20").Value = "Raw Influent Biological Oxygen Demand (mg/L)"
$wsConc.Range("C21").Value = "Raw Influent Chemical Oxygen Demand (mg/L)"
$wsConc.Range("C22").Value = "Raw Influent Total Kjedhal Nitrogen (mg/L)"
$wsConc.Range("C23").Value = "Raw Influent Total Nitrogen (mg/L)"
$wsConc.Range("C24").Value = "Raw Influent Total Phosphorus (mg/L)"
$wsConc.Range("C25").Value = "Raw Influent Total Suspended Solids (mg/L)"
$wsConc.Range("C26").Value = "Raw Influent Volatile Suspended Solids (mg/L)"

# --- decision_var_bound sheet (sheet4) ---
$wsBound = $wb.Worksheets.Item("decision_var_bound")
$wsBound.Range("B6").Value = 0

# --- fuzzy_goal sheet (sheet5) ---
$wsGoal = $wb.Worksheets.Item("fuzzy_goal")

# --- Update selections on each sheet (restore decision_var_bound as active tab at the end) ---
[void]$wsConc.Activate()
[void]$wsConc.Range("D27").Select()

[void]$wsGoal.Activate()
[void]$wsGoal.Range("F9").Select()

[void]$wsBound.Activate()
[void]$wsBound.Range("H11").Select()
